$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edits (user stories backlog) ---
# Row 4: "weet wat kan verwachten" -> "weet wat ik kan verwachten"
$ws.Range("F4").Value = "weet wat ik kan verwachten"

# Row 6: "zodat ik mijn intresses kan vergroten" -> "mijn intresses kan vergroten"
$ws.Range("F6").Value = "mijn intresses kan vergroten"

# Row 8: "tentoonstellingen kunnen bekijken" -> "info voor scholen kunnen vinden"
#        "weet wat ik kan zien in het museum" -> "mischien een school reis naar het museum kan regelen"
$ws.Range("D8").Value = "info voor scholen kunnen vinden"
$ws.Range("F8").Value = "mischien een school reis naar het museum kan regelen"

# --- Re-fit columns that now hold longer text ---
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(6).AutoFit()

# --- Update selection/active cell to reflect where the author was last working ---
$ws.Range("F6").Select()
